$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "_old" columns become "_FV2410", "_new" columns become "_FV2504".
#    Column K ("diff") is left untouched.
$headersFV2410 = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)
$headersFV2504 = @(
  "Segmentname_FV2504",
  "Segmentgruppe_FV2504",
  "Segment_FV2504",
  "Datenelement_FV2504",
  "Segment ID_FV2504",
  "Code_FV2504",
  "Qualifier_FV2504",
  "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504",
  "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headersFV2410[$i]
}
for ($i = 0; $i -lt $headersFV2504.Length; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $headersFV2504[$i]
}

# 2. Turn the used range into an Excel Table ("Table1") so the header row
#    gets AutoFilter + structured referencing.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true | Out-Null
